$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 11:52"

# 2. Update province/region data rows with new figures.
#    Row 17/18 also swap their region label (La Rioja now ranks above Albacete).

# Row 4: Madrid
$ws.Range("B4").Value = 57997
$ws.Range("C4").Value = 32277
$ws.Range("D4").Value = 18260
$ws.Range("E4").Value = 7460

# Row 5: Cataluña
$ws.Range("B5").Value = 43112
$ws.Range("C5").Value = 14446
$ws.Range("D5").Value = 24514
$ws.Range("E5").Value = 4152

# Row 6: Castilla-La Mancha
$ws.Range("B6").Value = 17045
$ws.Range("C6").Value = 4242
$ws.Range("D6").Value = 10728
$ws.Range("E6").Value = 2075

# Row 7: Castilla y Leon
$ws.Range("B7").Value = 16259
$ws.Range("C7").Value = 5429
$ws.Range("D7").Value = 9309
$ws.Range("E7").Value = 1521

# Row 8: Pais Vasco
$ws.Range("B8").Value = 12810
$ws.Range("C8").Value = 7277
$ws.Range("D8").Value = 4430
$ws.Range("E8").Value = 1103

# Row 9: Andalucia
$ws.Range("B9").Value = 11447
$ws.Range("C9").Value = 3411
$ws.Range("D9").Value = 7019
$ws.Range("E9").Value = 1017

# Row 10: Galicia
$ws.Range("B10").Value = 8468
$ws.Range("C10").Value = 1588
$ws.Range("D10").Value = 6520
$ws.Range("E10").Value = 360

# Row 14: Aragon
$ws.Range("B14").Value = 4938
$ws.Range("C14").Value = 1182
$ws.Range("D14").Value = 3113
$ws.Range("E14").Value = 643

# Row 15: Navarra
$ws.Range("B15").Value = 4781
$ws.Range("C15").Value = 1190
$ws.Range("D15").Value = 3194
$ws.Range("E15").Value = 397

# Row 17: now La Rioja (was Albacete)
$ws.Range("A17").Value = "La Rioja"
$ws.Range("B17").Value = 3763
$ws.Range("C17").Value = 1866
$ws.Range("D17").Value = 1603
$ws.Range("E17").Value = 294

# Row 18: now Albacete (was La Rioja)
$ws.Range("A18").Value = "Albacete"
$ws.Range("B18").Value = 3754
$ws.Range("C18").Value = 4178
$ws.Range("D18").Value = 10597
$ws.Range("E18").Value = 373

# Row 22: Extremadura
$ws.Range("B22").Value = 3196
$ws.Range("C22").Value = 921
$ws.Range("D22").Value = 1878
$ws.Range("E22").Value = 397

# Row 28: Asturias
$ws.Range("B28").Value = 2365
$ws.Range("C28").Value = 616
$ws.Range("D28").Value = 1547
$ws.Range("E28").Value = 202

# Row 32: Cantabria
$ws.Range("B32").Value = 2123
$ws.Range("C32").Value = 646
$ws.Range("D32").Value = 1312
$ws.Range("E32").Value = 165

# Row 33: Gran Canaria
$ws.Range("B33").Value = 2085
$ws.Range("C33").Value = 878
$ws.Range("D33").Value = 1087
$ws.Range("E33").Value = 120

# Row 36: Murcia
$ws.Range("B36").Value = 1681
$ws.Range("C36").Value = 702
$ws.Range("D36").Value = 859
$ws.Range("E36").Value = 120

# Row 58: Ceuta
$ws.Range("C58").Value = 73
$ws.Range("D58").Value = 34

# Row 59: Melilla
$ws.Range("B59").Value = 105
$ws.Range("D59").Value = 57
